$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two rows whose A-column sequence is no longer present in this
# pass (GGAAGACCTGATACC and TCTTTAATCCAGATA). Delete the higher row first so
# the lower row number stays valid.
$ws.Rows(29).Delete()
$ws.Rows(20).Delete()

# Update the B-column (quantification) values for the remaining 30 rows.
$ws.Range("B1").Value = 0
$ws.Range("B2").Value = 0.0337723741979
$ws.Range("B3").Value = 0.00667779632721202
$ws.Range("B4").Value = 0.01794043774668
$ws.Range("B5").Value = 0.06759040216289
$ws.Range("B6").Value = 0.01387347391786
$ws.Range("B7").Value = 0.04750593824228
$ws.Range("B8").Value = 0
$ws.Range("B9").Value = 0
$ws.Range("B10").Value = 0
$ws.Range("B11").Value = 4.377028148818949
$ws.Range("B12").Value = 0.01240233163834
$ws.Range("B13").Value = 0.1131221719457
$ws.Range("B14").Value = 0.0204081632653
$ws.Range("B17").Value = 0
$ws.Range("B18").Value = 0
$ws.Range("B19").Value = 0.01710132535271
$ws.Range("B20").Value = 0
$ws.Range("B21").Value = 0.009463423866754991
$ws.Range("B22").Value = 0.01688903901368
$ws.Range("B23").Value = 0.0262559075792
$ws.Range("B24").Value = 0.01785395465095
$ws.Range("B25").Value = 0.1360544217687
$ws.Range("B26").Value = 0.02965599051008
$ws.Range("B27").Value = 0
$ws.Range("B28").Value = 0.02645502645502
$ws.Range("B29").Value = 0
$ws.Range("B30").Value = 0.006619009796134498
